# Applies the odds updates for the FlashScore 2024-10-31 weekly matches sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (AS Roma vs Torino)
$ws.Range("I2").Value = 4.5
$ws.Range("L2").Value = 4.75
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 8.5
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.67
$ws.Range("AA2").Value = 17
$ws.Range("AC2").Value = 8.5
$ws.Range("AD2").Value = 6.5
$ws.Range("AG2").Value = 401
$ws.Range("AH2").Value = 11
$ws.Range("AO2").Value = 10
$ws.Range("AP2").Value = 23
$ws.Range("AT2").Value = 2.63
$ws.Range("AV2").Value = 67
$ws.Range("AX2").Value = 23
$ws.Range("BA2").Value = 126
$ws.Range("BC2").Value = 126

# Row 3 (Como vs Lazio)
$ws.Range("N3").Value = 10
$ws.Range("BC3").Value = 151

# Row 4 (Estudiantes L.P. vs Ind. Rivadavia)
$ws.Range("G4").Value = 1.53

# Row 5 (Sarmiento Junin vs Independiente)
$ws.Range("H5").Value = 2.82
$ws.Range("I5").Value = 2.15

# Row 6 (Aurora vs Independiente)
$ws.Range("G6").Value = 1.49

# Row 7 (Grasshoppers vs Lugano)
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 11
$ws.Range("Q7").Value = 1.85
$ws.Range("R7").Value = 2

# Row 8 (Servette vs Luzern)
$ws.Range("G8").Value = 1.8
$ws.Range("J8").Value = 2.38
$ws.Range("K8").Value = 2.38
